$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.390.77'
$ws.Range("E2").Value = '  -0.72%  '
$ws.Range("D3").Value = '1.848.72'
$ws.Range("E3").Value = '  -0.33%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9992'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.62'
$ws.Range("E5").Value = '  -0.88%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6346'
$ws.Range("E6").Value = '  -1.02%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07574'
$ws.Range("E8").Value = '  -0.53%  '
$ws.Range("E9").Value = '  -1.40%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.88'
$ws.Range("E10").Value = '  +1.70%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07742'
$ws.Range("E11").Value = '  +0.81%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.008'
$ws.Range("E12").Value = '  -0.98%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6834'
$ws.Range("E13").Value = '  -1.10%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '83.11'
$ws.Range("E14").Value = '  -1.31%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000009974'
$ws.Range("E15").Value = '  +2.59%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.142'
$ws.Range("E16").Value = '  -2.65%  '
$ws.Range("D17").Value = '29.423.29'
$ws.Range("E17").Value = '  -0.78%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '231.00'
$ws.Range("E18").Value = '  -3.37%  '
$ws.Range("E19").Value = '  -1.37%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.0000'
$ws.Range("E20").Value = '  -0.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.551'
$ws.Range("E21").Value = '  -1.34%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("E23").Value = '  +230.87%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '16.59'
$ws.Range("E24").Value = '  +170.49%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '156.53'
$ws.Range("E25").Value = '  -0.47%  '
$ws.Range("E26").Value = '  -0.47%  '
$ws.Range("E27").Value = '  -1.58%  '
$ws.Range("E28").Value = '  -0.62%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.733'
$ws.Range("E29").Value = '  +173.97%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.05703'
$ws.Range("E31").Value = '  -3.18%  '
$ws.Range("E32").Value = '  -2.44%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.134'
$ws.Range("E33").Value = '  -0.40%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.027'
$ws.Range("E34").Value = '  -1.63%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.848'
$ws.Range("E35").Value = '  -3.34%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.156'
$ws.Range("E36").Value = '  -2.61%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.7170'
$ws.Range("E37").Value = '  -1.25%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.600'
$ws.Range("E38").Value = '  +0.04%  '
$ws.Range("D39").Value = '1.246.54'
$ws.Range("E39").Value = '  +2.23%  '
$ws.Range("E40").Value = '  -0.09%  '
$ws.Range("E41").Value = '  +1.38%  '
$ws.Range("E42").Value = '  +263.31%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9018'
$ws.Range("E43").Value = '  -1.50%  '
$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.000'
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.89'
$ws.Range("E45").Value = '  -0.13%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '66.26'
$ws.Range("E46").Value = '  -1.98%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.170'
$ws.Range("E47").Value = '  -0.21%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.051'
$ws.Range("E48").Value = '  -6.04%  '
$ws.Range("E49").Value = '  +1.22%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4026'
$ws.Range("E50").Value = '  -1.13%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1127'
$ws.Range("E51").Value = '  -0.39%  '
